# Update the TPM-derived statistics in the LR-pair worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ E=3; F=1; G=0.3915976666666667; H=1.174793; I=0.02606065131430495; J=0.02606065131430495; K=3; L=1; M=1.315861666666667; N=3.947585; O=0.2754050739440597; P=0.2754050739440597; Q=0.5152883583227779; R=4.637595224905001; S=0.007177235602246513; T=0.007177235602246513 }
    3  = @{ E=3; F=1; G=0.3915976666666667; H=1.174793; I=0.02606065131430495; J=0.02606065131430495; O=0.3040809095127364; P=0.3040809095127364; Q=0.568941415697889; R=5.120472741281; S=0.00792454655414814; T=0.00792454655414814 }
    4  = @{ E=3; F=1; G=0.3915976666666667; H=1.174793; I=0.02606065131430495; J=0.02606065131430495; M=2.009179666666667; N=6.027539; O=0.4205140165432039; P=0.4205140165432039; Q=0.786790069380778; R=7.081110624427001; S=0.0109588691579103; T=0.0109588691579103 }
    5  = @{ I=0.4187506438669658; J=0.4187506438669658; K=3; L=1; M=1.315861666666667; N=3.947585; O=0.2754050739440597; P=0.2754050739440597; Q=8.279813471368334; R=74.518321242315; S=0.1153260520383043; T=0.1153260520383043 }
    6  = @{ I=0.4187506438669658; J=0.4187506438669658; O=0.3040809095127364; P=0.3040809095127364; S=0.1273340766461109; T=0.1273340766461109 }
    7  = @{ I=0.4187506438669658; J=0.4187506438669658; M=2.009179666666667; N=6.027539; O=0.4205140165432039; P=0.4205140165432039; Q=12.64238733590233; R=113.781486023121; S=0.1760905151825506; T=0.1760905151825506 }
    8  = @{ G=8.342485333333334; H=25.027456; I=0.5551887048187292; J=0.5551887048187292; K=3; L=1; M=1.315861666666667; N=3.947585; O=0.2754050739440597; P=0.2754050739440597; Q=10.97755665486222; R=98.79800989376001; S=0.1529017863035088; T=0.1529017863035088 }
    9  = @{ G=8.342485333333334; H=25.027456; I=0.5551887048187292; J=0.5551887048187292; O=0.3040809095127364; P=0.3040809095127364; Q=12.12056613203911; R=109.085095188352; S=0.1688222863124773; T=0.1688222863124773 }
    10 = @{ G=8.342485333333334; H=25.027456; I=0.5551887048187292; J=0.5551887048187292; M=2.009179666666667; N=6.027539; O=0.4205140165432039; P=0.4205140165432039; Q=16.76155190119822; R=150.853967110784; S=0.233464632202743; T=0.233464632202743 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
